$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$p.Style = "Heading 3"
